$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C ("Resource") - this shifts existing URL/ContentType/... columns
# from C..O to D..P, and the drawing/hyperlink/conditional-formatting targets move with them.
$ws.Columns.Item(3).Insert()

# Give the new column roughly the same width as column B (Excel normally carries the
# preceding column's width onto a freshly inserted column).
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Populate the new "Resource" column: header in row 1, "test" in the data rows.
$ws.Range("C1").Value = "Resource"
$ws.Range("C2").Value = "test"
$ws.Range("C3").Value = "test"
$ws.Range("C4").Value = "test"

# Match formatting of the new cells to the rest of the sheet's look (bold/green header-ish
# style used elsewhere for this workbook, copied from existing cells with the same look).
$ws.Range("B4").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "Resource"

$ws.Range("B4").Copy($ws.Range("C2"))
$ws.Range("C2").Value = "test"

$ws.Range("E4").Copy($ws.Range("C3"))
$ws.Range("C3").Value = "test"

$ws.Range("B4").Copy($ws.Range("C4"))
$ws.Range("C4").Value = "test"

# The conditional formatting that used to live on G2 now belongs on H2 (shifted by the new
# column); keep the original rule (type/dxf) but repoint it instead of recreating it so the
# notContainsBlanks rule and its dxf linkage survive untouched.
$fc = $ws.Range("G2").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("H2"))

# The three hyperlinks (previously on C2/C3/C4) now belong on D2/D3/D4. Re-create them there,
# using foreach (not .Item(n), which this host mis-resolves for Hyperlinks) to read the
# original target addresses before removing the stale entries.
$links = @()
foreach ($h in $ws.Hyperlinks) {
    $links += $h.Address
}
$ws.Hyperlinks.Delete()

# Adding a hyperlink stamps the built-in "Hyperlink" look onto the cell, clobbering the
# rich/bold formatting those cells already carry - stash and restore each cell's formatting
# around the Add() call.
$ws.Range("D2").Copy($ws.Range("Z90"))
$ws.Range("D3").Copy($ws.Range("Z91"))
$ws.Range("D4").Copy($ws.Range("Z92"))

$ws.Hyperlinks.Add($ws.Range("D2"), $links[0])
$ws.Hyperlinks.Add($ws.Range("D3"), $links[1])
$ws.Hyperlinks.Add($ws.Range("D4"), $links[2])

$ws.Range("Z90").Copy($ws.Range("D2"))
$ws.Range("Z91").Copy($ws.Range("D3"))
$ws.Range("Z92").Copy($ws.Range("D4"))
$ws.Range("Z90:Z92").ClearContents()

# Restore the active selection to the new "Resource" column's last data row, with the view
# scrolled one column to the right (as the saved workbook shows).
$ws.Range("C4").Select()
